$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.023.08'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '1.831.81'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '324.16'
$ws.Range("E5").Value = '  -3.51%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.4642'
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").Value = '0.3866'
$ws.Range("E8").Value = '  -1.17%  '
$ws.Range("D9").Value = '0.07838'
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").Value = '0.9596'
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("D11").Value = '21.89'
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").Value = '1.826.91'
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("D13").Value = '5.685'
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("D14").Value = '6.890'
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("D15").Value = '0.06858'
$ws.Range("E15").Value = '  -0.10%  '
$ws.Range("D16").Value = '88.18'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '0.000009892'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").Value = '16.62'
$ws.Range("E19").Value = '  -2.99%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").Value = '28.037.45'
$ws.Range("E21").Value = '  -2.16%  '
$ws.Range("D22").Value = '5.298'
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("D23").Value = '10.96'
$ws.Range("E23").Value = '  -3.33%  '
$ws.Range("D24").Value = '2.090'
$ws.Range("E24").Value = '  -2.40%  '
$ws.Range("D25").Value = '2.032.50'
$ws.Range("E25").Value = '  -5.63%  '
$ws.Range("D26").Value = '154.88'
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").Value = '19.12'
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("D28").Value = '5.650'
$ws.Range("E28").Value = '  -6.56%  '
$ws.Range("D29").Value = '1.957'
$ws.Range("E29").Value = '  -3.43%  '
$ws.Range("D30").Value = '118.22'
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("D31").Value = '0.09241'
$ws.Range("E31").Value = '  -1.71%  '
$ws.Range("D32").Value = '0.9328'
$ws.Range("E32").Value = '  -4.45%  '
$ws.Range("D33").Value = '5.253'
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("D34").Value = '1.318'
$ws.Range("E34").Value = '  -2.14%  '
$ws.Range("D35").Value = '3.306'
$ws.Range("E35").Value = '  -4.98%  '
$ws.Range("D36").Value = '0.05848'
$ws.Range("E36").Value = '  -4.70%  '
$ws.Range("D37").Value = '0.02123'
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("D38").Value = '1.144'
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("D39").Value = '7.749'
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("D40").Value = '0.5587'
$ws.Range("E40").Value = '  -2.49%  '
$ws.Range("D41").Value = '9.865'
$ws.Range("E41").Value = '  -3.05%  '
$ws.Range("D42").Value = '0.1758'
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("D43").Value = '0.07192'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("D44").Value = '11.60'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").Value = '0.5260'
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("D46").Value = '1.150'
$ws.Range("E46").Value = '  -8.12%  '
$ws.Range("D47").Value = '2.096'
$ws.Range("E47").Value = '  -11.10%  '
$ws.Range("D48").Value = '1.822'
$ws.Range("E48").Value = '  -4.75%  '
$ws.Range("D49").Value = '112.79'
$ws.Range("E49").Value = '  -2.79%  '
$ws.Range("D50").Value = '1.001'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '1.024'
$ws.Range("E51").Value = '  +0.20%  '
